# Realestate Update resale numbers 2023-06-28 13:44
# Appends one new data row (row 85) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

# Columns A and D look like dates/plain numbers ("2023-06-28", "26"), so
# Excel's normal type-inference would silently turn them into a date serial
# or a number. Force them to be stored as literal text (matching the rest
# of the column), then drop the temporary "@" text format so the cell
# keeps the sheet's default (unstyled) formatting.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "2023-06-28"
$ws.Cells.Item($row, 1).ClearFormats()

$ws.Cells.Item($row, 2).Value = "13:44:28"
$ws.Cells.Item($row, 3).Value = "Wednesday"

$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "26"
$ws.Cells.Item($row, 4).ClearFormats()

$ws.Cells.Item($row, 5).Value = 122969
$ws.Cells.Item($row, 6).Value = 134425
$ws.Cells.Item($row, 7).Value = 163635
$ws.Cells.Item($row, 8).Value = 134053
$ws.Cells.Item($row, 9).Value = 177215
$ws.Cells.Item($row, 10).Value = 114862
$ws.Cells.Item($row, 11).Value = 203918
$ws.Cells.Item($row, 12).Value = 226446
$ws.Cells.Item($row, 13).Value = 176142
$ws.Cells.Item($row, 14).Value = 104437
$ws.Cells.Item($row, 15).Value = 39716
$ws.Cells.Item($row, 16).Value = 33763
$ws.Cells.Item($row, 17).Value = 52367
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35765
$ws.Cells.Item($row, 20).Value = -1
